# Logic Component Sequence Diagram: rename the deletePerson(...) call
# shown in the sequence diagram to deleteStudent(...), matching the
# DeveloperGuide.adoc diagram update.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The call is drawn in "TextBox 77", which holds two runs:
#   run 1: "deletePerson"
#   run 2: "(p)"
# Locate it robustly by scanning shapes for the text instead of a
# hard-coded index, then replace just the "deletePerson" substring so the
# existing run formatting (font color, dirty/err flags, etc.) is kept
# untouched.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $text = $tr.Text
        $pos = $text.IndexOf("deletePerson")
        if ($pos -ge 0) {
            $old = "deletePerson"
            $sub = $tr.Characters($pos + 1, $old.Length)
            $sub.Text = "deleteStudent"
        }
    }
}
